$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.921.11'
$ws.Range("E2").Value = '  +0.76%  '
$ws.Range("D3").Value = '3.803.02'
$ws.Range("E3").Value = '  -0.89%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '443.81'
$ws.Range("E5").Value = '  +5.23%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.56'
$ws.Range("E6").Value = '  +13.10%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.622'
$ws.Range("E7").Value = '  +3.26%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.734'
$ws.Range("E9").Value = '  +2.64%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.151'
$ws.Range("E10").Value = '  -7.04%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0000310'
$ws.Range("E11").Value = '  -10.08%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '43.42'
$ws.Range("E12").Value = '  +8.36%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.36'
$ws.Range("E13").Value = '  +4.73%  '
$ws.Range("D14").Value = '4.398.61'
$ws.Range("E14").Value = '  -0.40%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.73'
$ws.Range("E15").Value = '  -6.95%  '
$ws.Range("D17").Value = '3.755.99'
$ws.Range("E17").Value = '  -2.01%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '19.87'
$ws.Range("E18").Value = '  +2.30%  '
$ws.Range("E19").Value = '  +7.38%  '
$ws.Range("D20").Value = '66.911.72'
$ws.Range("E20").Value = '  +0.75%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '418.43'
$ws.Range("E21").Value = '  +4.25%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '14.57'
$ws.Range("E22").Value = '  +2.81%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.27'
$ws.Range("E23").Value = '  +10.47%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '85.65'
$ws.Range("E24").Value = '  +2.38%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '37.14'
$ws.Range("E25").Value = '  +0.22%  '
$ws.Range("E26").Value = '  +7.47%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '5.54'
$ws.Range("E27").Value = '  -4.06%  '
$ws.Range("B28").Value = 'Filecoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.73'
$ws.Range("E28").Value = '  +3.82%  '
$ws.Range("B29").Value = 'RenderToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.42'
$ws.Range("E29").Value = '  +27.00%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '732.33'
$ws.Range("E30").Value = '  +4.80%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '13.86'
$ws.Range("E31").Value = '  +13.50%  '
$ws.Range("E32").Value = '  +11.16%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.75'
$ws.Range("E33").Value = '  +0.09%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '44.21'
$ws.Range("E34").Value = '  +17.81%  '
$ws.Range("E35").Value = '  +7.03%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '56.48'
$ws.Range("E36").Value = '  +3.28%  '
$ws.Range("E37").Value = '  +0.06%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.50'
$ws.Range("E38").Value = '  +25.30%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0479'
$ws.Range("E39").Value = '  +6.39%  '
$ws.Range("B40").Value = 'TheGraph'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.344'
$ws.Range("E40").Value = '  +19.51%  '
$ws.Range("B41").Value = 'ThetaToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.87'
$ws.Range("E41").Value = '  -0.86%  '
$ws.Range("B42").Value = 'Stellar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.141'
$ws.Range("E42").Value = '  +4.92%  '
$ws.Range("B43").Value = 'PEPE'
$ws.Range("C43").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D43").Value = '0.0₃0673'
$ws.Range("E43").Value = '  -10.88%  '
$ws.Range("E44").Value = '  +0.11%  '
$ws.Range("B45").Value = 'Fetch.AI'
$ws.Range("C45").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.56'
$ws.Range("E45").Value = '  +28.40%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.26'
$ws.Range("E46").Value = '  +7.64%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.35'
$ws.Range("E47").Value = '  +1.70%  '
$ws.Range("B48").Value = 'ARBITRUM'
$ws.Range("C48").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.11'
$ws.Range("E48").Value = '  +3.37%  '
$ws.Range("B49").Value = 'Monero'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '144.86'
$ws.Range("E49").Value = '  +1.00%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.63'
$ws.Range("E50").Value = '  +4.86%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.84'
$ws.Range("E51").Value = '  +4.49%  '
